# Updates cryptos list data (price and 1h volume change columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = "23.413.74"
$dCell.Style = "Normal"
$ws.Range("E2").Value = "  -0.20%  "

$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = "1.629.94"
$dCell.Style = "Normal"
$ws.Range("E3").Value = "  -0.56%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("E5").Value = "  -0.06%  "

$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = "303.91"
$dCell.Style = "Normal"
$ws.Range("E6").Value = "  -0.95%  "

$dCell = $ws.Range("D7")
$dCell.NumberFormat = "@"
$dCell.Value = "0.3776"
$dCell.Style = "Normal"
$ws.Range("E7").Value = "  +0.33%  "

$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = "0.3646"
$dCell.Style = "Normal"
$ws.Range("E8").Value = "  +0.07%  "

$dCell = $ws.Range("D9")
$dCell.NumberFormat = "@"
$dCell.Value = "51.05"
$dCell.Style = "Normal"
$ws.Range("E9").Value = "  -2.36%  "

$dCell = $ws.Range("D10")
$dCell.NumberFormat = "@"
$dCell.Value = "0.08221"
$dCell.Style = "Normal"
$ws.Range("E10").Value = "  +0.82%  "

$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = "1.233"
$dCell.Style = "Normal"
$ws.Range("E11").Value = "  -2.48%  "

$ws.Range("E12").Value = "  +0.01%  "

$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = "22.35"
$dCell.Style = "Normal"
$ws.Range("E13").Value = "  -2.49%  "

$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"
$dCell.Value = "6.538"
$dCell.Style = "Normal"
$ws.Range("E14").Value = "  -1.45%  "

$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = "0.00001248"
$dCell.Style = "Normal"
$ws.Range("E15").Value = "  -2.16%  "

$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = "7.326"
$dCell.Style = "Normal"
$ws.Range("E16").Value = "  -0.57%  "

$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = "1.627.98"
$dCell.Style = "Normal"
$ws.Range("E17").Value = "  -0.77%  "

$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = "93.99"
$dCell.Style = "Normal"
$ws.Range("E18").Value = "  -0.63%  "

$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = "0.06979"
$dCell.Style = "Normal"
$ws.Range("E19").Value = "  +0.75%  "

$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = "17.67"
$dCell.Style = "Normal"
$ws.Range("E20").Value = "  -2.62%  "

$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = "6.499"
$dCell.Style = "Normal"
$ws.Range("E21").Value = "  -0.76%  "

$ws.Range("E22").Value = "  -0.05%  "

$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = "12.69"
$dCell.Style = "Normal"
$ws.Range("E23").Value = "  -0.77%  "

$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = "23.405.38"
$dCell.Style = "Normal"
$ws.Range("E24").Value = "  -0.22%  "

$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = "3.125"
$dCell.Style = "Normal"
$ws.Range("E25").Value = "  +0.73%  "

$dCell = $ws.Range("D26")
$dCell.NumberFormat = "@"
$dCell.Value = "2.452"
$dCell.Style = "Normal"
$ws.Range("E26").Value = "  +1.28%  "

$dCell = $ws.Range("D27")
$dCell.NumberFormat = "@"
$dCell.Value = "21.36"
$dCell.Style = "Normal"
$ws.Range("E27").Value = "  +0.57%  "

$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = "150.04"
$dCell.Style = "Normal"
$ws.Range("E28").Value = "  -0.46%  "

$dCell = $ws.Range("D29")
$dCell.NumberFormat = "@"
$dCell.Value = "5.287"
$dCell.Style = "Normal"
$ws.Range("E29").Value = "  -1.40%  "

$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = "133.77"
$dCell.Style = "Normal"
$ws.Range("E30").Value = "  -1.08%  "

$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = "1.810.18"
$dCell.Style = "Normal"
$ws.Range("E31").Value = "  -0.80%  "

$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = "2.261"
$dCell.Style = "Normal"
$ws.Range("E32").Value = "  -1.75%  "

$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value = "6.823"
$dCell.Style = "Normal"
$ws.Range("E33").Value = "  +0.41%  "

$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = "1.028"
$dCell.Style = "Normal"
$ws.Range("E34").Value = "  +6.43%  "

$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = "10.80"
$dCell.Style = "Normal"
$ws.Range("E35").Value = "  +4.67%  "

$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = "0.02782"
$dCell.Style = "Normal"
$ws.Range("E36").Value = "  -1.39%  "

$dCell = $ws.Range("D37")
$dCell.NumberFormat = "@"
$dCell.Value = "0.2520"
$dCell.Style = "Normal"
$ws.Range("E37").Value = "  -0.22%  "

$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = "0.08758"
$dCell.Style = "Normal"
$ws.Range("E38").Value = "  -0.89%  "

$dCell = $ws.Range("D39")
$dCell.NumberFormat = "@"
$dCell.Value = "0.07100"
$dCell.Style = "Normal"
$ws.Range("E39").Value = "  -2.75%  "

$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = "5.992"
$dCell.Style = "Normal"
$ws.Range("E40").Value = "  -2.09%  "

$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = "0.7026"
$dCell.Style = "Normal"
$ws.Range("E41").Value = "  -1.13%  "

$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = "1.345"
$dCell.Style = "Normal"
$ws.Range("E42").Value = "  -2.16%  "

$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = "16.17"
$dCell.Style = "Normal"
$ws.Range("E43").Value = "  +0.17%  "

$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = "12.11"
$dCell.Style = "Normal"
$ws.Range("E44").Value = "  -3.27%  "

$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = "0.6541"
$dCell.Style = "Normal"
$ws.Range("E45").Value = "  -0.02%  "

$ws.Range("E46").Value = "  -0.07%  "

$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = "2.294"
$dCell.Style = "Normal"
$ws.Range("E47").Value = "  -1.88%  "

$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = "3.976"
$dCell.Style = "Normal"
$ws.Range("E48").Value = "  -1.13%  "

$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = "0.08015"
$dCell.Style = "Normal"
$ws.Range("E49").Value = "  +0.70%  "

$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = "1.200"
$dCell.Style = "Normal"
$ws.Range("E50").Value = "  -0.06%  "

$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = "125.83"
$dCell.Style = "Normal"
$ws.Range("E51").Value = "  -2.35%  "
